# Auto-generated PowerShell COM-interop script to apply the diff
# to Bahamut_Profits workbook sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 75.44444
$ws.Range("I4").Value = 59.875
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 59.875
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 54.125
$ws.Range("N4").Value = -428
$ws.Range("H18").Value = 1630
$ws.Range("I18").Value = 1195
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 1195
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = -911
$ws.Range("N18").Value = -3068
$ws.Range("H132").Value = 265076.1
$ws.Range("I132").Value = 2107.5
$ws.Range("K132").Value = 6322.5
$ws.Range("M132").Value = -3792.5
$ws.Range("H135").Value = 1631.5
$ws.Range("I135").Value = 845.1177
$ws.Range("K135").Value = 7606.0593
$ws.Range("M135").Value = -5071.0593
$ws.Range("H137").Value = 1044.091
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 1164.3334
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 3493.0002
$ws.Range("M137").Value = -447
$ws.Range("N137").Value = -8593.0002
$ws.Range("H138").Value = 3256.05
$ws.Range("I138").Value = 883.3333
$ws.Range("J138").Value = 3674.7646
$ws.Range("K138").Value = 2649.9999
$ws.Range("L138").Value = 11024.2938
$ws.Range("M138").Value = 2490.0001
$ws.Range("N138").Value = -21304.2938

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6416.511
$ws.Range("I32").Value = 5658.6577
$ws.Range("J32").Value = 10530.571
$ws.Range("K32").Value = 5658.6577
$ws.Range("L32").Value = 10530.571
$ws.Range("M32").Value = -5371.6577
$ws.Range("N32").Value = -11104.571
$ws.Range("H61").Value = 2013.1177
$ws.Range("I61").Value = 2173.1428
$ws.Range("K61").Value = 2173.1428
$ws.Range("M61").Value = -1961.1428
$ws.Range("H110").Value = 709.4706
$ws.Range("I110").Value = 670.8461
$ws.Range("J110").Value = 835
$ws.Range("K110").Value = 670.8461
$ws.Range("L110").Value = 835
$ws.Range("M110").Value = 1374.1539
$ws.Range("N110").Value = -4925
$ws.Range("H136").Value = 2013.1177
$ws.Range("I136").Value = 2173.1428
$ws.Range("K136").Value = 6519.428400000001
$ws.Range("M136").Value = -3969.428400000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1166.5834
$ws.Range("I16").Value = 1166.5834
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1166.5834
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -879.5834
$ws.Range("H58").Value = 5374.087
$ws.Range("I58").Value = 778.8570999999999
$ws.Range("K58").Value = 778.8570999999999
$ws.Range("M58").Value = -575.8570999999999
$ws.Range("H64").Value = 29800
$ws.Range("J64").Value = 29800
$ws.Range("L64").Value = 29800
$ws.Range("N64").Value = -30296
$ws.Range("H67").Value = 29800
$ws.Range("J67").Value = 29800
$ws.Range("L67").Value = 29800
$ws.Range("N67").Value = -31516
$ws.Range("H93").Value = 8571
$ws.Range("I93").Value = 3336.375
$ws.Range("J93").Value = 50448
$ws.Range("K93").Value = 3336.375
$ws.Range("L93").Value = 50448
$ws.Range("M93").Value = -1464.375
$ws.Range("N93").Value = -54192
$ws.Range("H113").Value = 1166.5834
$ws.Range("I113").Value = 1166.5834
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1166.5834
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1003.4166
$ws.Range("H134").Value = 2755.7368
$ws.Range("I134").Value = 1625.9
$ws.Range("J134").Value = 4011.111
$ws.Range("K134").Value = 4877.700000000001
$ws.Range("L134").Value = 12033.333
$ws.Range("M134").Value = -2342.700000000001
$ws.Range("N134").Value = -17103.333
$ws.Range("H136").Value = 5374.087
$ws.Range("I136").Value = 778.8570999999999
$ws.Range("K136").Value = 2336.5713
$ws.Range("M136").Value = 213.4287000000004
$ws.Range("N16").ClearContents()
$ws.Range("N113").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4014.5625
$ws.Range("I3").Value = 1470
$ws.Range("J3").Value = 8255.5
$ws.Range("K3").Value = 4410
$ws.Range("L3").Value = 24766.5
$ws.Range("M3").Value = -4298
$ws.Range("N3").Value = -24990.5
$ws.Range("H5").Value = 1245.6
$ws.Range("I5").Value = 742.6
$ws.Range("J5").Value = 3257.6
$ws.Range("K5").Value = 2227.8
$ws.Range("L5").Value = 9772.799999999999
$ws.Range("M5").Value = -2115.8
$ws.Range("N5").Value = -9996.799999999999
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H113").Value = 575.38464
$ws.Range("I113").Value = 543.75
$ws.Range("J113").Value = 583.5484
$ws.Range("K113").Value = 1631.25
$ws.Range("L113").Value = 1750.6452
$ws.Range("M113").Value = 538.75
$ws.Range("N113").Value = -6090.6452
$ws.Range("H135").Value = 1245.6
$ws.Range("I135").Value = 742.6
$ws.Range("J135").Value = 3257.6
$ws.Range("K135").Value = 6683.400000000001
$ws.Range("L135").Value = 29318.4
$ws.Range("M135").Value = -4148.400000000001
$ws.Range("N135").Value = -34388.39999999999
$ws.Range("M51").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3075.0454
$ws.Range("I80").Value = 3167.0588
$ws.Range("J80").Value = 2762.2
$ws.Range("K80").Value = 3167.0588
$ws.Range("L80").Value = 2762.2
$ws.Range("M80").Value = -2169.0588
$ws.Range("N80").Value = -4758.2
$ws.Range("H83").Value = 3075.0454
$ws.Range("I83").Value = 3167.0588
$ws.Range("J83").Value = 2762.2
$ws.Range("K83").Value = 15835.294
$ws.Range("L83").Value = 13811
$ws.Range("M83").Value = -10843.294
$ws.Range("N83").Value = -23795
$ws.Range("H97").Value = 2089.4119
$ws.Range("I97").Value = 2092.4167
$ws.Range("J97").Value = 2082.2
$ws.Range("K97").Value = 2092.4167
$ws.Range("L97").Value = 2082.2
$ws.Range("M97").Value = -1596.4167
$ws.Range("N97").Value = -3074.2

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 18000
$ws.Range("J106").Value = 18000
$ws.Range("L106").Value = 18000
$ws.Range("N106").Value = -20524

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2911.1
$ws.Range("I122").Value = 2901.2222
$ws.Range("K122").Value = 8703.6666
$ws.Range("M122").Value = -6253.6666
